$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ui")

# Copy the last row (row 12, which carries the "bottom border" style) down to the
# new row 13 first, so the new row inherits that distinctive formatting.
$ws.Range("A12:B12").Copy()
$ws.Range("A13:B13").PasteSpecial(-4122)

# Row 12 is no longer the last row, so it should go back to the regular style
# used by rows 3-11 (copy formatting only from row 11).
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)

# Match the row height used by the other data rows.
$ws.Rows.Item(13).RowHeight = 20

# Fill in the new row's content: a continuation of the shared ROW()-2 formula
# and the new menu-button label.
$ws.Range("A13").Formula = "=ROW()-2"
$ws.Range("B13").Value = "メニュー"

$wb.Application.CutCopyMode = $false
